# Update "想去人数" (F) / "最低票价" (G) figures to the newer scrape snapshot.
# Sheet order (per xl/workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value2  = 27052
$ws1.Cells.Item(4, 6).Value2  = 654
$ws1.Cells.Item(5, 6).Value2  = 195
$ws1.Cells.Item(6, 6).Value2  = 577
$ws1.Cells.Item(8, 6).Value2  = 377
$ws1.Cells.Item(9, 6).Value2  = 482
$ws1.Cells.Item(10, 6).Value2 = 197
$ws1.Cells.Item(12, 6).Value2 = 316
$ws1.Cells.Item(13, 6).Value2 = 101
$ws1.Cells.Item(14, 6).Value2 = 506
$ws1.Cells.Item(16, 6).Value2 = 1638
$ws1.Cells.Item(17, 6).Value2 = 264
$ws1.Cells.Item(18, 6).Value2 = 886
$ws1.Cells.Item(19, 6).Value2 = 194
$ws1.Cells.Item(20, 6).Value2 = 463
$ws1.Cells.Item(22, 6).Value2 = 109

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value2  = 4528
$ws2.Cells.Item(3, 6).Value2  = 248
$ws2.Cells.Item(3, 7).Value2  = "不可售"
$ws2.Cells.Item(9, 6).Value2  = 5
$ws2.Cells.Item(11, 6).Value2 = 458

# ---- Sheet 3: 本地生活 ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value2 = 5194
$ws3.Cells.Item(3, 6).Value2 = 279

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value2  = 5194
$ws4.Cells.Item(4, 6).Value2  = 279
$ws4.Cells.Item(5, 6).Value2  = 27052
$ws4.Cells.Item(6, 6).Value2  = 4528
$ws4.Cells.Item(7, 6).Value2  = 248
$ws4.Cells.Item(7, 7).Value2  = "不可售"
$ws4.Cells.Item(8, 6).Value2  = 654
$ws4.Cells.Item(11, 6).Value2 = 195
$ws4.Cells.Item(15, 6).Value2 = 5
$ws4.Cells.Item(17, 6).Value2 = 458
$ws4.Cells.Item(18, 6).Value2 = 577
$ws4.Cells.Item(22, 6).Value2 = 377
$ws4.Cells.Item(23, 6).Value2 = 482
$ws4.Cells.Item(24, 6).Value2 = 197
$ws4.Cells.Item(27, 6).Value2 = 316
$ws4.Cells.Item(28, 6).Value2 = 101
$ws4.Cells.Item(31, 6).Value2 = 506
$ws4.Cells.Item(34, 6).Value2 = 1638
$ws4.Cells.Item(35, 6).Value2 = 265
$ws4.Cells.Item(36, 6).Value2 = 887
$ws4.Cells.Item(38, 6).Value2 = 194
$ws4.Cells.Item(39, 6).Value2 = 463
$ws4.Cells.Item(41, 6).Value2 = 109

$wb.Save()
